$wb = $excel.ActiveWorkbook

# ALC row 51
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 6250
$ws.Range("J51").Value = 6500
$ws.Range("L51").Value = 6500
$ws.Range("N51").Value = -7468

# ALC row 76
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3299.6667
$ws.Range("I76").Value = 3000
$ws.Range("K76").Value = 3000
$ws.Range("M76").Value = -2685

# ALC row 79
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 3299.6667
$ws.Range("I79").Value = 3000
$ws.Range("K79").Value = 3000
$ws.Range("M79").Value = -1908

# ALC row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2029.8182
$ws.Range("J112").Value = 2029.8182
$ws.Range("L112").Value = 6089.4546
$ws.Range("N112").Value = -8305.454600000001

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 964.7692
$ws.Range("I132").Value = 912.36365
$ws.Range("J132").Value = 1253
$ws.Range("K132").Value = 2737.09095
$ws.Range("L132").Value = 3759
$ws.Range("M132").Value = -207.0909499999998
$ws.Range("N132").Value = -8819

# ALC row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 2035.25
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 2035.25
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 18317.25
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -23387.25

# ALC row 136
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H136").Value = 68149.89
$ws.Range("J136").Value = 68149.89
$ws.Range("L136").Value = 68149.89
$ws.Range("N136").Value = -78349.89

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1914.5
$ws.Range("I137").Value = 1647.875
$ws.Range("K137").Value = 4943.625
$ws.Range("M137").Value = -2393.625

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1827.4667
$ws.Range("I138").Value = 1304.28
$ws.Range("J138").Value = 2201.1714
$ws.Range("K138").Value = 3912.84
$ws.Range("L138").Value = 6603.514200000001
$ws.Range("M138").Value = 1227.16
$ws.Range("N138").Value = -16883.5142

# ALC row 139
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H139").Value = 69920
$ws.Range("J139").Value = 69920
$ws.Range("L139").Value = 69920
$ws.Range("N139").Value = -80200

# ALC row 140
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H140").Value = 72034.89999999999
$ws.Range("J140").Value = 72034.89999999999
$ws.Range("L140").Value = 72034.89999999999
$ws.Range("N140").Value = -82394.89999999999

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5524.656
$ws.Range("I32").Value = 3310.34
$ws.Range("K32").Value = 3310.34
$ws.Range("M32").Value = -3023.34

# ARM row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("N63").ClearContents()

# ARM row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("N66").ClearContents()

# ARM row 88
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 5650
$ws.Range("J88").Value = 5650
$ws.Range("L88").Value = 5650
$ws.Range("N88").Value = -6462

# ARM row 91
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 5650
$ws.Range("J91").Value = 5650
$ws.Range("L91").Value = 5650
$ws.Range("N91").Value = -8458

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1030.0869
$ws.Range("I122").Value = 842.5333000000001
$ws.Range("J122").Value = 1381.75
$ws.Range("K122").Value = 2527.5999
$ws.Range("L122").Value = 4145.25
$ws.Range("M122").Value = -77.59990000000016
$ws.Range("N122").Value = -9045.25

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1484.3513
$ws.Range("I132").Value = 1147
$ws.Range("K132").Value = 3441
$ws.Range("M132").Value = -911

# BSM row 64
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 773.75
$ws.Range("I64").Value = 698.3333
$ws.Range("J64").Value = 1000
$ws.Range("K64").Value = 698.3333
$ws.Range("L64").Value = 1000
$ws.Range("M64").Value = -473.3333
$ws.Range("N64").Value = -1450

# BSM row 67
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H67").Value = 773.75
$ws.Range("I67").Value = 698.3333
$ws.Range("J67").Value = 1000
$ws.Range("K67").Value = 698.3333
$ws.Range("L67").Value = 1000
$ws.Range("M67").Value = 81.66669999999999
$ws.Range("N67").Value = -2560

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2096.7917
$ws.Range("I105").Value = 1883.4762
$ws.Range("J105").Value = 3590
$ws.Range("K105").Value = 1883.4762
$ws.Range("L105").Value = 3590
$ws.Range("M105").Value = -136.4762000000001
$ws.Range("N105").Value = -7084

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 961.1429000000001
$ws.Range("I107").Value = 631.1667
$ws.Range("K107").Value = 631.1667
$ws.Range("M107").Value = 1288.8333

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5647.9614
$ws.Range("I134").Value = 6234.476
$ws.Range("K134").Value = 18703.428
$ws.Range("M134").Value = -16168.428

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2439.2666
$ws.Range("I31").Value = 2277.8
$ws.Range("J31").Value = 2520
$ws.Range("K31").Value = 2277.8
$ws.Range("L31").Value = 2520
$ws.Range("M31").Value = -1982.8
$ws.Range("N31").Value = -3110

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2439.2666
$ws.Range("I34").Value = 2277.8
$ws.Range("J34").Value = 2520
$ws.Range("K34").Value = 2277.8
$ws.Range("L34").Value = 2520
$ws.Range("M34").Value = -2075.8
$ws.Range("N34").Value = -2924

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4350247
$ws.Range("I58").Value = 8696454
$ws.Range("K58").Value = 8696454
$ws.Range("M58").Value = -8696251

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3681.7
$ws.Range("I132").Value = 1603
$ws.Range("J132").Value = 4572.5713
$ws.Range("K132").Value = 4809
$ws.Range("L132").Value = 13717.7139
$ws.Range("M132").Value = -2279
$ws.Range("N132").Value = -18777.7139

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 4350247
$ws.Range("I136").Value = 8696454
$ws.Range("K136").Value = 26089362
$ws.Range("M136").Value = -26086812

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 17812.268
$ws.Range("J131").Value = 19158.5
$ws.Range("L131").Value = 57475.5
$ws.Range("N131").Value = -67555.5

# GSM row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1336.1
$ws.Range("I113").Value = 1084.75
$ws.Range("K113").Value = 1084.75
$ws.Range("M113").Value = 1085.25

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1454.619
$ws.Range("I122").Value = 1455.1428
$ws.Range("J122").Value = 1453.5714
$ws.Range("K122").Value = 4365.428400000001
$ws.Range("L122").Value = 4360.7142
$ws.Range("M122").Value = -1915.428400000001
$ws.Range("N122").Value = -9260.7142

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2461429.2
$ws.Range("J126").Value = 113554.22
$ws.Range("L126").Value = 340662.66
$ws.Range("N126").Value = -345602.66

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2750864
$ws.Range("I132").Value = 4809900
$ws.Range("J132").Value = 5482.5
$ws.Range("K132").Value = 14429700
$ws.Range("L132").Value = 16447.5
$ws.Range("M132").Value = -14427170
$ws.Range("N132").Value = -21507.5

# GSM row 134
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H134").Value = 30171.143
$ws.Range("J134").Value = 30171.143
$ws.Range("L134").Value = 90513.429
$ws.Range("N134").Value = -95583.429

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2271.15
$ws.Range("I7").Value = 2162.389
$ws.Range("J7").Value = 3250
$ws.Range("K7").Value = 2162.389
$ws.Range("L7").Value = 3250
$ws.Range("M7").Value = -2050.389
$ws.Range("N7").Value = -3474

# LTW row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 7371.6665
$ws.Range("I16").Value = 8168.125
$ws.Range("K16").Value = 8168.125
$ws.Range("M16").Value = -7998.125

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 12298.866
$ws.Range("I40").Value = 13418.9
$ws.Range("J40").Value = 10058.8
$ws.Range("K40").Value = 13418.9
$ws.Range("L40").Value = 10058.8
$ws.Range("M40").Value = -13282.9
$ws.Range("N40").Value = -10330.8

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2395.4375
$ws.Range("I61").Value = 2391
$ws.Range("J61").Value = 2399.875
$ws.Range("K61").Value = 2391
$ws.Range("L61").Value = 2399.875
$ws.Range("M61").Value = -2189
$ws.Range("N61").Value = -2803.875

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 2395.4375
$ws.Range("I113").Value = 2391
$ws.Range("J113").Value = 2399.875
$ws.Range("K113").Value = 2391
$ws.Range("L113").Value = 2399.875
$ws.Range("M113").Value = -221
$ws.Range("N113").Value = -6739.875

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2271.15
$ws.Range("I126").Value = 2162.389
$ws.Range("J126").Value = 3250
$ws.Range("K126").Value = 6487.167
$ws.Range("L126").Value = 9750
$ws.Range("M126").Value = -4017.167
$ws.Range("N126").Value = -14690

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2256.6
$ws.Range("I132").Value = 1431.0714
$ws.Range("K132").Value = 4293.2142
$ws.Range("M132").Value = -1763.2142

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 7748.75
$ws.Range("I136").Value = 5497.5
$ws.Range("K136").Value = 16492.5
$ws.Range("M136").Value = -13942.5

# WVR row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 853.4167
$ws.Range("I107").Value = 471.33334
$ws.Range("K107").Value = 1414.00002
$ws.Range("M107").Value = 505.9999800000001

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 36495.227
$ws.Range("I122").Value = 44116.777
$ws.Range("K122").Value = 132350.331
$ws.Range("M122").Value = -129900.331

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1753.4667
$ws.Range("I126").Value = 1650.9
$ws.Range("K126").Value = 4952.700000000001
$ws.Range("M126").Value = -2482.700000000001

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1115.5471
$ws.Range("I132").Value = 780.59576
$ws.Range("J132").Value = 3739.3333
$ws.Range("K132").Value = 2341.78728
$ws.Range("L132").Value = 11217.9999
$ws.Range("M132").Value = 188.21272
$ws.Range("N132").Value = -16277.9999

# WVR row 133
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value = 69143.336
$ws.Range("J133").Value = 80000
$ws.Range("L133").Value = 80000
$ws.Range("N133").Value = -90120

# WVR row 140
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H140").Value = 62300
$ws.Range("J140").Value = 62300
$ws.Range("L140").Value = 62300
$ws.Range("N140").Value = -72660
